$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 427 (existing rows 427:527 shift down to 428:528)
$ws.Rows.Item(427).Insert()

# Populate the new row 427 with the new record's data
$ws.Range("A427").Value = 4
$ws.Range("B427").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C427").Value = "Los Lagos"
$ws.Range("D427").Value = 45204
$ws.Range("E427").Value = 10
$ws.Range("F427").Value = 100112003
$ws.Range("G427").Value = "Ajo"
$ws.Range("H427").Value = "Chino"
$ws.Range("I427").Value = "Primera"
$ws.Range("J427").Value = 120
$ws.Range("K427").Value = 25000
$ws.Range("L427").Value = 25000
$ws.Range("M427").Value = 25000
$ws.Range("N427").Value = "$/caja 10 kilos"
$ws.Range("O427").Value = "China"
$ws.Range("P427").Value = 2500
$ws.Range("Q427").Value = 10
$ws.Range("R427").Value = "Hortaliza"
